$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99; this shifts the existing rows 99-158
# down to 100-159 (and the sheet dimension grows from R158 to R159).
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new data record.
$ws.Cells.Item(99, 1).Value = 7
$ws.Cells.Item(99, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(99, 3).Value = "Ñuble"
$ws.Cells.Item(99, 4).Value = 44455
$ws.Cells.Item(99, 5).Value = 16
$ws.Cells.Item(99, 6).Value = 100112023
$ws.Cells.Item(99, 7).Value = "Brócoli"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 600
$ws.Cells.Item(99, 11).Value = 750
$ws.Cells.Item(99, 12).Value = 800
$ws.Cells.Item(99, 13).Value = 775
$ws.Cells.Item(99, 14).Value = "$/unidad"
$ws.Cells.Item(99, 15).Value = "Región del Maule"
$ws.Cells.Item(99, 16).Value = 775
$ws.Cells.Item(99, 17).Value = 1
$ws.Cells.Item(99, 18).Value = "Hortaliza"
